# Status and queries updated for next session
# Adds the newest status row (2020-03-26 / "Implementation" / "SurveyEngineCore
# Implementation" / new comment) to the bottom of the "Status(Summary)" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status(Summary)")
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row - this keeps the table's ref/autoFilter ranges,
# the worksheet dimension and the totals in sync automatically.
$newListRow = $lo.ListRows.Add()
$newRowIndex = $newListRow.Range.Row
$lastRowIndex = $newRowIndex - 1

# Copy the formatting (number format / fill / alignment) from the previous
# last row down onto the freshly inserted row before filling in values.
$ws.Range($ws.Cells.Item($lastRowIndex, 1), $ws.Cells.Item($lastRowIndex, 4)).Copy() | Out-Null
$newRowRange = $ws.Range($ws.Cells.Item($newRowIndex, 1), $ws.Cells.Item($newRowIndex, 4))
$newRowRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new status entry.
$ws.Cells.Item($newRowIndex, 1).Value2 = 43916
$ws.Cells.Item($newRowIndex, 2).Value2 = "Implementation"
$ws.Cells.Item($newRowIndex, 3).Value2 = "SurveyEngineCore Implementation"
$ws.Cells.Item($newRowIndex, 4).Value2 = "1) Discussion on upcoming queries`n2) Wrote a few flowcharts for surveyEngine like constructor,setTimestampFor,addRenderedItem and will be changed after some reviews`n3) Implemented ``setTimeStampFor`` function and tests"

# Match the row height Excel settles on for the new wrapped comment, and move
# the active selection one row below the new data, same as the saved file.
$ws.Rows.Item($newRowIndex).RowHeight = 75
$ws.Cells.Item($newRowIndex + 1, 4).Select() | Out-Null
